$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "ActiveCell Before:" $ws.Range("E1").Value
$ws.Range("E1").Value = "MTTR(horas)"
Write-Host "ActiveCell After:" $ws.Range("E1").Value
